# Weekly update: insert a new "Sandia" price record for the week of
# 2021-11-05 (serial 44505) at the top of the historical block, pushing
# the existing rows 185-214 down by one (to 186-215).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(185).Insert()

$ws.Range("A185").Value = 9
$ws.Range("B185").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C185").Value = "Metropolitana"
$ws.Range("D185").Value = 44505
$ws.Range("E185").Value = 13
$ws.Range("F185").Value = 100112028
$ws.Range("G185").Value = "Sandia"
$ws.Range("H185").Value = "Sin especificar"
$ws.Range("I185").Value = "Primera"
$ws.Range("J185").Value = 160
$ws.Range("K185").Value = 800
$ws.Range("L185").Value = 1000
$ws.Range("M185").Value = 900
$ws.Range("N185").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O185").Value = "Perú"
$ws.Range("P185").Value = 900
$ws.Range("Q185").Value = 1
$ws.Range("R185").Value = "Hortaliza"
